# Natmi following Dr Hou advice:
# recompute the Fgf2-Fgfr2 LR-pair sheet with the revised "3 nearest
# neighbours" statistics, overwriting rows 2-7 and appending the new
# sCs-vs-{ECs,FAPs,sCs} rows 8-10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ row=2; A="ECs"; B="Fgf2"; C="Fgfr2"; D="ECs"; E=3.0; F=1.0; G=0.747119; H=2.241357; I=0.03096954854571248; J=0.03096954854571248; K=2.0; L=0.6666666666666666; M=0.09434; N=0.28302; O=0.05191071108246543; P=0.05191071108246543; Q=0.07048320645999999; R=0.63434885814; S=0.001607651286910868; T=0.001607651286910868 },
    @{ row=3; A="ECs"; B="Fgf2"; C="Fgfr2"; D="FAPs"; E=3.0; F=1.0; G=0.747119; H=2.241357; I=0.03096954854571248; J=0.03096954854571248; K=3.0; L=1.0; M=0.9431116666666667; N=2.829335; O=0.5189484550226392; P=0.5189484550226391; Q=0.7046166452883333; R=6.341549807594999; S=0.01607159937054611; T=0.01607159937054611 },
    @{ row=4; A="ECs"; B="Fgf2"; C="Fgfr2"; D="sCs"; E=3.0; F=1.0; G=0.747119; H=2.241357; I=0.03096954854571248; J=0.03096954854571248; K=3.0; L=1.0; M=0.7798996666666667; N=2.339699; O=0.4291408338948954; P=0.4291408338948954; Q=0.5826778590603333; R=5.244100731543; S=0.0132902978882555; T=0.0132902978882555 },
    @{ row=5; A="FAPs"; B="Fgf2"; C="Fgfr2"; D="ECs"; E=3.0; F=1.0; G=19.74619233333334; H=59.23857700000001; I=0.8185184181638298; J=0.8185184181638298; K=2.0; L=0.6666666666666666; M=0.09434; N=0.28302; O=0.05191071108246543; P=0.05191071108246543; Q=1.862855784726667; R=16.76570206254; S=0.04248987312097919; T=0.04248987312097919 },
    @{ row=6; A="FAPs"; B="Fgf2"; C="Fgfr2"; D="FAPs"; E=3.0; F=1.0; G=19.74619233333334; H=59.23857700000001; I=0.8185184181638298; J=0.8185184181638298; K=3.0; L=1.0; M=0.9431116666666667; N=2.829335; O=0.5189484550226392; P=0.5189484550226391; Q=18.62286436181056; R=167.605779256295; S=0.424768868513694; T=0.4247688685136939 },
    @{ row=7; A="FAPs"; B="Fgf2"; C="Fgfr2"; D="sCs"; E=3.0; F=1.0; G=19.74619233333334; H=59.23857700000001; I=0.8185184181638298; J=0.8185184181638298; K=3.0; L=1.0; M=0.7798996666666667; N=2.339699; O=0.4291408338948954; P=0.4291408338948954; Q=15.40004881870256; R=138.600439368323; S=0.3512596765291566; T=0.3512596765291566 },
    @{ row=8; A="sCs"; B="Fgf2"; C="Fgfr2"; D="ECs"; E=3.0; F=1.0; G=3.630999; H=10.892997; I=0.1505120332904577; J=0.1505120332904577; K=2.0; L=0.6666666666666666; M=0.09434; N=0.28302; O=0.05191071108246543; P=0.05191071108246543; Q=0.3425484456599999; R=3.08293601094; S=0.007813186674575368; T=0.007813186674575368 },
    @{ row=9; A="sCs"; B="Fgf2"; C="Fgfr2"; D="FAPs"; E=3.0; F=1.0; G=3.630999; H=10.892997; I=0.1505120332904577; J=0.1505120332904577; K=3.0; L=1.0; M=0.9431116666666667; N=2.829335; O=0.5189484550226392; P=0.5189484550226391; Q=3.424437518555; R=30.819937666995; S=0.07810798713839905; T=0.07810798713839905 },
    @{ row=10; A="sCs"; B="Fgf2"; C="Fgfr2"; D="sCs"; E=3.0; F=1.0; G=3.630999; H=10.892997; I=0.1505120332904577; J=0.1505120332904577; K=3.0; L=1.0; M=0.7798996666666667; N=2.339699; O=0.4291408338948954; P=0.4291408338948954; Q=2.831814909767; R=25.486334187903; S=0.06459085947748326; T=0.06459085947748328 }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in $rows) {
    foreach ($c in $cols) {
        $ws.Range($c + $r.row).Value = $r.$c
    }
}

